# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)
#
# Turns the "Data" sheet (Suriname name/title + a Micro/SMEs/MSMEs detail
# row) into a "Summary" sheet: rename the tab, drop the Micro/SMEs/MSMEs
# detail row (B5:D5), and register the "title_" (bold+underline) named
# cell style alongside the existing name/title/source/HyperLink styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "Data" -> "Summary"
$ws.Name = "Summary"

# Remove the Micro / SMEs / MSMEs breakdown row entirely (not just clear
# it) so the used range shrinks back down to A1:A3 and the now-orphaned
# shared strings ("Micro", "SMEs", "MSMEs") drop out of the workbook.
$ws.Range("B5:D5").Delete()

# Register the "title_" cell style (bold + underlined Calibri 11) in the
# workbook's style gallery, next to the other named styles (name, title,
# source, HyperLink).
$titleStyle = $wb.Styles.Add("title_")
$titleStyle.Font.Bold = $true
$titleStyle.Font.Underline = $true
